$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.70913233333333
$ws.Range("H2").Value = 32.127397
$ws.Range("I2").Value = 0.007451729107954897
$ws.Range("J2").Value = 0.007451729107954897
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 1707.931377204008
$ws.Range("R2").Value = 15371.38239483607
$ws.Range("S2").Value = 0.002223063230449385
$ws.Range("T2").Value = 0.002223063230449385

$ws.Range("G3").Value = 10.70913233333333
$ws.Range("H3").Value = 32.127397
$ws.Range("I3").Value = 0.007451729107954897
$ws.Range("J3").Value = 0.007451729107954897
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 1847.951779614103
$ws.Range("R3").Value = 16631.56601652693
$ws.Range("S3").Value = 0.002405315405370011
$ws.Range("T3").Value = 0.002405315405370011

$ws.Range("G4").Value = 10.70913233333333
$ws.Range("H4").Value = 32.127397
$ws.Range("I4").Value = 0.007451729107954897
$ws.Range("J4").Value = 0.007451729107954897
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 796.6277839573833
$ws.Range("R4").Value = 7169.650055616449
$ws.Range("S4").Value = 0.001036899935505137
$ws.Range("T4").Value = 0.001036899935505137

$ws.Range("G5").Value = 10.70913233333333
$ws.Range("H5").Value = 32.127397
$ws.Range("I5").Value = 0.007451729107954897
$ws.Range("J5").Value = 0.007451729107954897
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 625.5698350962969
$ws.Range("R5").Value = 5630.128515866671
$ws.Range("S5").Value = 0.0008142489312173049
$ws.Range("T5").Value = 0.000814248931217305

$ws.Range("G6").Value = 10.70913233333333
$ws.Range("H6").Value = 32.127397
$ws.Range("I6").Value = 0.007451729107954897
$ws.Range("J6").Value = 0.007451729107954897
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 746.9214568932508
$ws.Range("R6").Value = 6722.293112039257
$ws.Range("S6").Value = 0.000972201605413058
$ws.Range("T6").Value = 0.0009722016054130582

$ws.Range("I7").Value = 0.03290895798513831
$ws.Range("J7").Value = 0.03290895798513832
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 7542.711378746259
$ws.Range("R7").Value = 67884.40240871633
$ws.Range("S7").Value = 0.009817680351673819
$ws.Range("T7").Value = 0.009817680351673821

$ws.Range("I8").Value = 0.03290895798513831
$ws.Range("J8").Value = 0.03290895798513832
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.01062255786134607
$ws.Range("T8").Value = 0.01062255786134607

$ws.Range("I9").Value = 0.03290895798513831
$ws.Range("J9").Value = 0.03290895798513832
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 3518.13517268911
$ws.Range("R9").Value = 31663.21655420199
$ws.Range("S9").Value = 0.004579245423173496
$ws.Range("T9").Value = 0.004579245423173498

$ws.Range("I10").Value = 0.03290895798513831
$ws.Range("J10").Value = 0.03290895798513832
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 2762.694553399289
$ws.Range("R10").Value = 24864.25098059361
$ws.Range("S10").Value = 0.00359595517747265
$ws.Range("T10").Value = 0.003595955177472652

$ws.Range("I11").Value = 0.03290895798513831
$ws.Range("J11").Value = 0.03290895798513832
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 3298.617876065587
$ws.Range("R11").Value = 29687.56088459028
$ws.Range("S11").Value = 0.004293519171472275
$ws.Range("T11").Value = 0.004293519171472276

$ws.Range("G12").Value = 411.37678
$ws.Range("H12").Value = 1234.13034
$ws.Range("I12").Value = 0.2862480573072345
$ws.Range("J12").Value = 0.2862480573072345
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 65607.86518887446
$ws.Range("R12").Value = 590470.78669987
$ws.Range("S12").Value = 0.08539595599469194
$ws.Range("T12").Value = 0.08539595599469194

$ws.Range("G13").Value = 411.37678
$ws.Range("H13").Value = 1234.13034
$ws.Range("I13").Value = 0.2862480573072345
$ws.Range("J13").Value = 0.2862480573072345
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 70986.55885749965
$ws.Range("R13").Value = 638879.0297174968
$ws.Range("S13").Value = 0.0923969258709795
$ws.Range("T13").Value = 0.0923969258709795

$ws.Range("G14").Value = 411.37678
$ws.Range("H14").Value = 1234.13034
$ws.Range("I14").Value = 0.2862480573072345
$ws.Range("J14").Value = 0.2862480573072345
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 30601.37482874109
$ws.Range("R14").Value = 275412.3734586698
$ws.Range("S14").Value = 0.03983110333996037
$ws.Range("T14").Value = 0.03983110333996037

$ws.Range("G15").Value = 411.37678
$ws.Range("H15").Value = 1234.13034
$ws.Range("I15").Value = 0.2862480573072345
$ws.Range("J15").Value = 0.2862480573072345
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 24030.41594938851
$ws.Range("R15").Value = 216273.7435444966
$ws.Range("S15").Value = 0.03127826727848039
$ws.Range("T15").Value = 0.0312782672784804

$ws.Range("G16").Value = 411.37678
$ws.Range("H16").Value = 1234.13034
$ws.Range("I16").Value = 0.2862480573072345
$ws.Range("J16").Value = 0.2862480573072345
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 28691.97375526448
$ws.Range("R16").Value = 258227.7637973803
$ws.Range("S16").Value = 0.03734580482312224
$ws.Range("T16").Value = 0.03734580482312225

$ws.Range("G17").Value = 173.2560603333334
$ws.Range("H17").Value = 519.768181
$ws.Range("I17").Value = 0.12055666021578
$ws.Range("J17").Value = 0.12055666021578
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 27631.50669200346
$ws.Range("R17").Value = 248683.5602280311
$ws.Range("S17").Value = 0.03596548863073658
$ws.Range("T17").Value = 0.03596548863073658

$ws.Range("G18").Value = 173.2560603333334
$ws.Range("H18").Value = 519.768181
$ws.Range("I18").Value = 0.12055666021578
$ws.Range("J18").Value = 0.12055666021578
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 29896.80536726132
$ws.Range("R18").Value = 269071.2483053519
$ws.Range("S18").Value = 0.03891402758152016
$ws.Range("T18").Value = 0.03891402758152016

$ws.Range("G19").Value = 173.2560603333334
$ws.Range("H19").Value = 519.768181
$ws.Range("I19").Value = 0.12055666021578
$ws.Range("J19").Value = 0.12055666021578
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 12888.12082104224
$ws.Range("R19").Value = 115993.0873893802
$ws.Range("S19").Value = 0.01677532709408491
$ws.Range("T19").Value = 0.01677532709408492

$ws.Range("G20").Value = 173.2560603333334
$ws.Range("H20").Value = 519.768181
$ws.Range("I20").Value = 0.12055666021578
$ws.Range("J20").Value = 0.12055666021578
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 10120.6859452844
$ws.Range("R20").Value = 91086.17350755959
$ws.Range("S20").Value = 0.01317320185821505
$ws.Range("T20").Value = 0.01317320185821506

$ws.Range("G21").Value = 173.2560603333334
$ws.Range("H21").Value = 519.768181
$ws.Range("I21").Value = 0.12055666021578
$ws.Range("J21").Value = 0.12055666021578
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 12083.9546073177
$ws.Range("R21").Value = 108755.5914658593
$ws.Range("S21").Value = 0.01572861505122326
$ws.Range("T21").Value = 0.01572861505122326

$ws.Range("G22").Value = 794.4973246666667
$ws.Range("H22").Value = 2383.491974
$ws.Range("I22").Value = 0.5528345953838922
$ws.Range("J22").Value = 0.5528345953838923
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 126709.3231894423
$ws.Range("R22").Value = 1140383.90870498
$ws.Range("S22").Value = 0.164926320282674
$ws.Range("T22").Value = 0.1649263202826741

$ws.Range("G23").Value = 794.4973246666667
$ws.Range("H23").Value = 2383.491974
$ws.Range("I23").Value = 0.5528345953838922
$ws.Range("J23").Value = 0.5528345953838923
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 137097.264215001
$ws.Range("R23").Value = 1233875.377935009
$ws.Range("S23").Value = 0.1784473844438121
$ws.Range("T23").Value = 0.1784473844438122

$ws.Range("G24").Value = 794.4973246666667
$ws.Range("H24").Value = 2383.491974
$ws.Range("I24").Value = 0.5528345953838922
$ws.Range("J24").Value = 0.5528345953838923
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 59100.8331402581
$ws.Range("R24").Value = 531907.498262323
$ws.Range("S24").Value = 0.0769263278353242
$ws.Range("T24").Value = 0.07692632783532423

$ws.Range("G25").Value = 794.4973246666667
$ws.Range("H25").Value = 2383.491974
$ws.Range("I25").Value = 0.5528345953838922
$ws.Range("J25").Value = 0.5528345953838923
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 46410.25480926846
$ws.Range("R25").Value = 417692.2932834161
$ws.Range("S25").Value = 0.06040812433059934
$ws.Range("T25").Value = 0.06040812433059935

$ws.Range("G26").Value = 794.4973246666667
$ws.Range("H26").Value = 2383.491974
$ws.Range("I26").Value = 0.5528345953838922
$ws.Range("J26").Value = 0.5528345953838923
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 55413.18201762347
$ws.Range("R26").Value = 498718.6381586112
$ws.Range("S26").Value = 0.07212643849148247
$ws.Range("T26").Value = 0.07212643849148249
